$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28: updated M/M % Delta figures for Durable Orders ---
$ws.Range("F28").Value = 0.0292251268148207
$ws.Range("G28").Value = -0.02746655290430811

# --- Row 29: updated Y/Y % Delta + 5yr,5yr Forward (T5YIFR) refresh ---
$ws.Range("G29").Value = 0.03397968857635882
$ws.Range("N29").Value = 45971
$ws.Range("Q29").Value = 2.2
$ws.Range("R29").Value = 2.2
$ws.Range("S29").Value = 2.19
$ws.Range("T29").Value = 2.21
$ws.Range("U29").Value = 2.2

# --- Row 30: 10yr TIPS (T10YIE) refresh ---
$ws.Range("N30").Value = 45971
$ws.Range("Q30").Value = 2.29
$ws.Range("R30").Value = 2.28
$ws.Range("S30").Value = 2.28
$ws.Range("T30").Value = 2.3
$ws.Range("U30").Value = 2.29

# --- Row 39: Nominal Broad US Dollar Index (DTWEXBGS) refresh ---
# Date cell needs to pick up the "recently updated" highlighted style (same
# look as N47:N52, style index 48) instead of the plain style it had before.
$ws.Range("N47").Copy() | Out-Null
$ws.Range("N39").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("N39").Value = 45968
$ws.Range("Q39").Value = 121.7835
$ws.Range("R39").Value = 122.0788
$ws.Range("S39").Value = 122.2295
$ws.Range("T39").Value = 122.2066
$ws.Range("U39").Value = 121.8422

# --- Row 47: FFR (DFF) refresh ---
$ws.Range("N47").Value = 45968
$ws.Range("U47").Value = 3.87

# --- Row 48: 2y UST (DGS2) refresh ---
$ws.Range("N48").Value = 45968
$ws.Range("Q48").Value = 3.55
$ws.Range("R48").Value = 3.57
$ws.Range("S48").Value = 3.63
$ws.Range("T48").Value = 3.58
$ws.Range("U48").Value = 3.6

# --- Row 49: 5y UST (DGS5) refresh ---
$ws.Range("N49").Value = 45968
$ws.Range("Q49").Value = 3.67
$ws.Range("R49").Value = 3.69
$ws.Range("S49").Value = 3.76
$ws.Range("T49").Value = 3.69
$ws.Range("U49").Value = 3.72

# --- Row 50: 10y UST (DGS10) refresh ---
$ws.Range("N50").Value = 45968
$ws.Range("Q50").Value = 4.11
$ws.Range("R50").Value = 4.11
$ws.Range("S50").Value = 4.17
$ws.Range("T50").Value = 4.1
$ws.Range("U50").Value = 4.13

# --- Row 52: BAA (DBAA) refresh ---
$ws.Range("N52").Value = 45968
$ws.Range("Q52").Value = 5.86
$ws.Range("R52").Value = 5.83
$ws.Range("S52").Value = 5.87
$ws.Range("T52").Value = 5.82
$ws.Range("U52").Value = 5.84
